$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (character-width offset to land on exact XML width values)
$offset = 0.8359375
$ws.Columns.Item(3).ColumnWidth = 8 - $offset
$ws.Columns.Item(10).ColumnWidth = 7 - $offset
$ws.Columns.Item(15).ColumnWidth = 8 - $offset
$ws.Columns.Item(28).ColumnWidth = 7 - $offset
$ws.Columns.Item(34).ColumnWidth = 8 - $offset

# Update data rows 2-5 with new values
$ws.Range("A2").Value = 45055.50694444445
$ws.Range("B2").Value = 13.283
$ws.Range("C2").Value = 8.791
$ws.Range("D2").Value = 3.507
$ws.Range("E2").Value = 28.828
$ws.Range("F2").Value = 21.492
$ws.Range("G2").Value = 10.237
$ws.Range("H2").Value = 30.575
$ws.Range("I2").Value = 16.288
$ws.Range("J2").Value = 6.473
$ws.Range("K2").Value = 9.559
$ws.Range("L2").Value = 11.331
$ws.Range("M2").Value = 12.141
$ws.Range("N2").Value = 3.376
$ws.Range("O2").Value = 10.527
$ws.Range("P2").Value = 14.419
$ws.Range("Q2").Value = 9.413
$ws.Range("R2").Value = 2.798
$ws.Range("S2").Value = 1.636
$ws.Range("T2").Value = 152.85
$ws.Range("U2").Value = 29.125
$ws.Range("V2").Value = 9.717
$ws.Range("W2").Value = 18.774
$ws.Range("X2").Value = 9.568
$ws.Range("Y2").Value = 2.805
$ws.Range("Z2").Value = 16.47
$ws.Range("AA2").Value = 8.583
$ws.Range("AB2").Value = 7.902
$ws.Range("AC2").Value = 9.369
$ws.Range("AD2").Value = 11.722
$ws.Range("AE2").Value = 3.066
$ws.Range("AF2").Value = 27.695
$ws.Range("AG2").Value = 5.181
$ws.Range("AH2").Value = 12.147
$ws.Range("A3").Value = 45055.51388888889
$ws.Range("B3").Value = 5.631
$ws.Range("C3").Value = 3.763
$ws.Range("D3").Value = 1.4
$ws.Range("E3").Value = 12.496
$ws.Range("F3").Value = 9.102
$ws.Range("G3").Value = 4.309
$ws.Range("H3").Value = 19.292
$ws.Range("I3").Value = 6.981
$ws.Range("J3").Value = 2.805
$ws.Range("K3").Value = 3.875
$ws.Range("L3").Value = 4.939
$ws.Range("M3").Value = 5.359
$ws.Range("N3").Value = 1.454
$ws.Range("O3").Value = 4.512
$ws.Range("P3").Value = 6.196
$ws.Range("Q3").Value = 4.246
$ws.Range("R3").Value = 1.3
$ws.Range("S3").Value = 0.701
$ws.Range("T3").Value = 61.369
$ws.Range("U3").Value = 12.783
$ws.Range("V3").Value = 4.164
$ws.Range("W3").Value = 8.136
$ws.Range("X3").Value = 4.113
$ws.Range("Y3").Value = 1.2
$ws.Range("Z3").Value = 9.505
$ws.Range("AA3").Value = 3.678
$ws.Range("AB3").Value = 3.488
$ws.Range("AC3").Value = 4.109
$ws.Range("AD3").Value = 5.065
$ws.Range("AE3").Value = 1.196
$ws.Range("AF3").Value = 18.211
$ws.Range("AG3").Value = 2.141
$ws.Range("AH3").Value = 5.208
$ws.Range("A4").Value = 45055.52083333334
$ws.Range("B4").Value = 16.703
$ws.Range("C4").Value = 12.285
$ws.Range("D4").Value = 1.314
$ws.Range("E4").Value = 36.572
$ws.Range("F4").Value = 29.397
$ws.Range("G4").Value = 13.063
$ws.Range("H4").Value = 47.963
$ws.Range("I4").Value = 20.36
$ws.Range("J4").Value = 8.908
$ws.Range("K4").Value = 13.088
$ws.Range("L4").Value = 14.639
$ws.Range("M4").Value = 15.604
$ws.Range("N4").Value = 4.225
$ws.Range("O4").Value = 13.159
$ws.Range("P4").Value = 18.6
$ws.Range("Q4").Value = 11.273
$ws.Range("R4").Value = 0.971
$ws.Range("S4").Value = 0.823
$ws.Range("T4").Value = 192.952
$ws.Range("U4").Value = 36.705
$ws.Range("V4").Value = 12.146
$ws.Range("W4").Value = 24.482
$ws.Range("X4").Value = 12.823
$ws.Range("Y4").Value = 2.128
$ws.Range("Z4").Value = 23.822
$ws.Range("AA4").Value = 10.728
$ws.Range("AB4").Value = 9.598
$ws.Range("AC4").Value = 11.291
$ws.Range("AD4").Value = 15.315
$ws.Range("AE4").Value = 0.745
$ws.Range("AF4").Value = 43.603
$ws.Range("AG4").Value = 6.755
$ws.Range("AH4").Value = 15.185
$ws.Range("A5").Value = 45055.52777777778
$ws.Range("B5").Value = 2.79
$ws.Range("C5").Value = 1.92
$ws.Range("D5").Value = 0.64
$ws.Range("E5").Value = 6.27
$ws.Range("F5").Value = 4.55
$ws.Range("G5").Value = 2.14
$ws.Range("H5").Value = 14.62
$ws.Range("I5").Value = 3.49
$ws.Range("J5").Value = 1.49
$ws.Range("K5").Value = 1.9
$ws.Range("L5").Value = 2.5
$ws.Range("M5").Value = 2.74
$ws.Range("N5").Value = 0.74
$ws.Range("O5").Value = 2.26
$ws.Range("P5").Value = 3.21
$ws.Range("Q5").Value = 2.16
$ws.Range("R5").Value = 0.65
$ws.Range("S5").Value = 0.32
$ws.Range("T5").Value = 27.1
$ws.Range("U5").Value = 6.71
$ws.Range("V5").Value = 2.08
$ws.Range("W5").Value = 4.36
$ws.Range("X5").Value = 2.14
$ws.Range("Y5").Value = 0.6
$ws.Range("Z5").Value = 6.7
$ws.Range("AA5").Value = 1.84
$ws.Range("AB5").Value = 1.77
$ws.Range("AC5").Value = 2.06
$ws.Range("AD5").Value = 2.57
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 13.98
$ws.Range("AG5").Value = 1.05
$ws.Range("AH5").Value = 2.61

# Remove row 6 (dataset now has one fewer row)
$ws.Rows.Item(6).Delete()
